$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $searchText, $replaceText) {
    $p = $d.Paragraphs($paraIndex).Range
    $rng = $d.Range($p.Start, $p.End)
    $found = $rng.Find.Execute(
        $searchText, $true, $false, $false, $false, $false, $true, 1, $false,
        "", 0
    )
    if (-not $found) {
        throw "Text not found in paragraph $paraIndex : $searchText"
    }
    $rng.Text = $replaceText
}

$quote = [char]34

# --- Background paragraph (paragraph 4: "Background: " + long intro text) ---
Replace-InParagraph 4 "Background: " "Background:"

$bgSearch = "Introducing " + $quote + "EcoWorx" + $quote + ", a new and innovative company founded in 2022 with a mission to promote eco-friendly solutions for modern-day problems. The company's goal is to provide sustainable alternatives to traditional products that often harm the environment. EcoWorx sells a variety of products, ranging from biodegradable packaging materials, eco-friendly cleaning products, and even reusable water bottles made from recycled plastic. The company was started by a group of environmentally conscious entrepreneurs who wanted to make a positive impact on the planet while also promoting responsible consumption. Their vision is to create a world where sustainability is at the forefront of every decision-making process, and consumers have easy access to affordable and eco-friendly products. "
$bgReplace = "EcoWorx is an innovative company founded in 2022 with a mission to promote eco-friendly solutions. The company seeks to provide sustainable alternatives to traditional products and services that often harm the environment. "
Replace-InParagraph 4 $bgSearch $bgReplace

# --- Objective paragraph (paragraph 6: long objective text) ---
$objSearch = "The primary objective of EcoWorx is to promote eco-friendly solutions and increase the availability of sustainable products. Company growth and expansion is a priority, and EcoWorx is committed to minimizing its carbon footprint by sourcing materials locally and ensuring that all products are manufactured in an environmentally responsible way. "
$objReplace = "To implement a new line of sewer waste treatment products in order to increase company growth and profitability. "
Replace-InParagraph 6 $objSearch $objReplace

# --- Target Audience paragraph (paragraph 8: "Target Audience: " + long text) ---
Replace-InParagraph 8 "Target Audience: " "Target Audience:"

$taSearch = "The target audience of EcoWorx is anyone who is interested in living an eco-friendly lifestyle and wants to contribute to environmental sustainability. EcoWorx appeals to all ages, incomes, and genders, as they are dedicated to making sustainable products accessible to everyone."
$taReplace = "The primary target audience is individuals and businesses who are looking to make a positive change in their lifestyle and take a more sustainability-driven approach to their everyday activities.  Additionally, EcoWorx also hopes to reach out to government organizations and larger businesses who are looking for more eco-sustainable products for their operations. "
Replace-InParagraph 8 $taSearch $taReplace

# --- Brand Guidelines paragraph (paragraph 10: "Brand Guidelines:" + long text) ---
Replace-InParagraph 10 "Brand Guidelines:" "Brand Guidelines: "

$bgdSearch = "EcoWorx is an environmentally conscious company that believes in the power of small changes to make a big difference. They embrace their core values of sustainability, responsibility, and positivity in order to create an inspiring, inviting, and eco-friendly brand. Their logo and visuals should express these values and be visually appealing to their target audience. Communication should be informative, respectful, and inspiring, with a focus on promoting eco-friendly solutions and encouraging sustainability."
$bgdReplace = "EcoWorx seeks to create a brand identity that speaks to the company's commitment to providing eco-friendly products and services. This should be communicated through the company's logo, slogan, color palette, messaging, and overall visual identity. The company should incorporate messaging that resonates with its target audience and evokes a feeling of environmental responsibility. Eco-friendly materials and resources should be used to create branded assets whenever possible."
Replace-InParagraph 10 $bgdSearch $bgdReplace
